$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data update: insert two new rows (269-270) for the latest week,
# pushing the previously-existing rows 269-275 down to 271-277.
$ws.Rows("269:270").Insert()

# Row 269: new week entry - Primera
$ws.Range("A269").Value = 4
$ws.Range("B269").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C269").Value = "Los Lagos"
$ws.Range("D269").Value = 44747
$ws.Range("E269").Value = 10
$ws.Range("F269").Value = "Fruta"
$ws.Range("G269").Value = 100104
$ws.Range("H269").Value = "Frutos de pepita"
$ws.Range("I269").Value = 100104005
$ws.Range("J269").Value = "Pera"
$ws.Range("K269").Value = "Packham's Triumph"
$ws.Range("L269").Value = "Primera"
$ws.Range("M269").Value = 400
$ws.Range("N269").Value = 15000
$ws.Range("O269").Value = 15000
$ws.Range("P269").Value = 15000
$ws.Range("Q269").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R269").Value = "Región de O'Higgins"
$ws.Range("S269").Value = 1000
$ws.Range("T269").Value = 15

# Row 270: new week entry - Segunda
$ws.Range("A270").Value = 4
$ws.Range("B270").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C270").Value = "Los Lagos"
$ws.Range("D270").Value = 44747
$ws.Range("E270").Value = 10
$ws.Range("F270").Value = "Fruta"
$ws.Range("G270").Value = 100104
$ws.Range("H270").Value = "Frutos de pepita"
$ws.Range("I270").Value = 100104005
$ws.Range("J270").Value = "Pera"
$ws.Range("K270").Value = "Packham's Triumph"
$ws.Range("L270").Value = "Segunda"
$ws.Range("M270").Value = 200
$ws.Range("N270").Value = 12000
$ws.Range("O270").Value = 12000
$ws.Range("P270").Value = 12000
$ws.Range("Q270").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R270").Value = "Región de O'Higgins"
$ws.Range("S270").Value = 800
$ws.Range("T270").Value = 15
